$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 7)
$ws.Range("A7").Value = "184. Department Highest Salary"
$ws.Range("B7").Value = "Medium"
$ws.Range("C7").Value = "Data Manipulation"
$ws.Range("D7").Value = "Merge the dataframes, group by department, then find the employees with the highest salary within each group using max function and boolean indexing. Handle empty table scenarios. Use the lambda function to find the highest salary in each group."
$ws.Range("E7").Value = "https://leetcode.com/problems/department-highest-salary/solutions/3861495/pandas-very-simple-step-by-step-process-detailed/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata "

# Match formatting used by the other "Medium" rows (B5/B6 fill) and hyperlink cells (E2..E6)
$ws.Range("B7").Interior.Color = 49407

# Add the hyperlink for the new row, then restore the shared Hyperlink cell style
$ws.Hyperlinks.Add($ws.Range("E7"), "https://leetcode.com/problems/department-highest-salary/solutions/3861495/pandas-very-simple-step-by-step-process-detailed/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata ")
$ws.Range("E7").Style = "Hyperlink"

# Expand the table to include the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E7"))

# Update the active selection to match the saved workbook state
$ws.Range("E16").Select()
